$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency data (prices, links, names, and 1h volume % changes)
# All target cells are text cells (inlineStr in the original), so force text
# number format before assigning values to avoid Excel auto-converting numeric-
# looking strings (e.g. "312.52") into real numbers.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '42.247.58'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +0.06%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.281.87'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  -0.50%  '
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  +0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '312.52'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.49%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '102.08'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.15%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.622'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -0.96%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.47%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.597'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -1.30%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '38.68'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -2.11%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0895'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -1.09%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.22'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -2.14%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  +1.66%  '
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +1.89%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '14.98'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  -1.72%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.629.86'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  -0.46%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.282.59'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -0.27%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '42.441.27'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  +0.48%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.21'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  -2.32%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.29'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +7.14%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '72.85'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -0.69%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.48'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -1.66%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '263.28'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -4.44%  '
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -4.29%  '
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  +0.35%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.63'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  -1.73%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.93'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +15.13%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.34'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -0.95%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '22.33'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -1.60%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '35.83'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -4.32%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '164.86'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.33%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0861'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -1.44%  '
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -2.63%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.61'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -1.74%  '
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -5.29%  '
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  -1.74%  '
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -3.85%  '
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  -0.40%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -3.33%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +5.16%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '98.17'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +2.54%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '68.77'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  -1.29%  '
$ws.Range('B44').NumberFormat = '@'
$ws.Range('B44').Value = 'FirstDigitalUSD'
$ws.Range('C44').NumberFormat = '@'
$ws.Range('C44').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.00'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.32%  '
$ws.Range('B45').NumberFormat = '@'
$ws.Range('B45').Value = 'Algorand'
$ws.Range('C45').NumberFormat = '@'
$ws.Range('C45').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.225'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  +0.63%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '11.89'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  -0.39%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.700.42'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  +6.38%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '78.52'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.65%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '109.65'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -2.00%  '
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.13%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '8.62'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -3.55%  '
